# Apply cryptocurrency price/volume updates from the latest GitHub Actions scrape run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.653.35"
$ws.Range("E2").Value = "  -0.48%  "

$ws.Range("D3").Value = "3.808.04"
$ws.Range("E3").Value = "  +2.05%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Formula = "'611.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.77%  "

$ws.Range("D6").Formula = "'176.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.61%  "

$ws.Range("D7").Value = "3.802.78"
$ws.Range("E7").Value = "  +1.96%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("E9").Value = "  -1.83%  "

$ws.Range("D10").Formula = "'0.166"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.62%  "

$ws.Range("D11").Formula = "'6.45"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.34%  "

$ws.Range("D13").Formula = "'40.04"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.32%  "

$ws.Range("E14").Value = "  -2.32%  "

$ws.Range("D15").Value = "4.442.84"
$ws.Range("E15").Value = "  +2.08%  "

$ws.Range("D16").Value = "3.813.45"
$ws.Range("E16").Value = "  +2.13%  "

$ws.Range("D17").Value = "69.693.72"
$ws.Range("E17").Value = "  -0.46%  "

$ws.Range("D18").Formula = "'7.50"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.73%  "

$ws.Range("E19").Value = "  -3.47%  "

$ws.Range("D20").Formula = "'16.62"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.88%  "

$ws.Range("D21").Formula = "'505.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.03%  "

$ws.Range("D22").Formula = "'9.52"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.12%  "

$ws.Range("D23").Formula = "'0.736"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.91%  "

$ws.Range("D24").Formula = "'85.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.79%  "

$ws.Range("D25").Formula = "'2.44"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.87%  "

$ws.Range("D26").Formula = "'0.0000143"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.55%  "

$ws.Range("E27").Value = "  -3.72%  "

$ws.Range("D28").Formula = "'10.42"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.91%  "

$ws.Range("E29").Value = "  +0.11%  "

$ws.Range("E30").Value = "  +1.58%  "

$ws.Range("E31").Value = "  +1.48%  "

$ws.Range("D32").Formula = "'7.95"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.43%  "

$ws.Range("D33").Formula = "'31.59"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.41%  "

$ws.Range("E34").Value = "  -1.89%  "

$ws.Range("E35").Value = "  -0.02%  "

$ws.Range("D36").Formula = "'1.04"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.05%  "

$ws.Range("D37").Formula = "'6.08"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.07%  "

$ws.Range("E38").Value = "  +3.52%  "

$ws.Range("D39").Formula = "'484.22"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +14.16%  "

$ws.Range("E40").Value = "  +0.03%  "

$ws.Range("D41").Formula = "'3.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.81%  "

$ws.Range("D42").Formula = "'2.04"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.08%  "

$ws.Range("D43").Formula = "'49.70"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.43%  "

$ws.Range("D44").Formula = "'43.65"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.93%  "

$ws.Range("D45").Formula = "'8.51"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.14%  "

$ws.Range("D46").Value = "2.917.51"
$ws.Range("E46").Value = "  -2.62%  "

$ws.Range("E47").Value = "  -0.93%  "

$ws.Range("D48").Formula = "'139.77"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.14%  "

$ws.Range("D50").Formula = "'26.76"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.12%  "

$ws.Range("E51").Value = "  -4.20%  "
